# fix: updated fake data WP5 and updated categories WP5
$wb = $excel.ActiveWorkbook

$wsVariables = $wb.Worksheets.Item("Variables")
$wsCategories = $wb.Worksheets.Item("Categories")

# --- Categories sheet: correct the category labels ---
# "asthma_current_CHICOS" -> "asthma_current_ISAAC" (rows 24-25)
$wsCategories.Range("A24").Value = "asthma_current_ISAAC"
$wsCategories.Range("A25").Value = "asthma_current_ISAAC"

# "pets_pregn" -> "pets_preg" (rows 46-47)
$wsCategories.Range("A46").Value = "pets_preg"
$wsCategories.Range("A47").Value = "pets_preg"

# --- Categories sheet: widen column A to fit the longer labels ---
$wsCategories.Columns.Item(1).ColumnWidth = 54.67

# --- Make "Categories" the active/selected sheet & cell ---
$wsCategories.Activate()
$wsCategories.Range("A46").Select()
